# Insert a new weekly price record for "Perejil" (Terminal La Palmera de La Serena)
# as row 15, pushing the existing data rows (old rows 15-120) down by one row
# (to new rows 16-121). This mirrors a new week of data being prepended to the
# historical series while keeping the rest of the rows unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 (entire row insert shifts rows 15:120 -> 16:121
# and copies formatting, e.g. the date style on column D, from the row below).
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Range("A15").Value2 = 8
$ws.Range("B15").Value2 = 'Terminal La Palmera de La Serena'
$ws.Range("C15").Value2 = 'Coquimbo'
$ws.Range("D15").Value2 = 44602
$ws.Range("E15").Value2 = 4
$ws.Range("F15").Value2 = 100112044
$ws.Range("G15").Value2 = 'Perejil'
$ws.Range("H15").Value2 = 'Sin especificar'
$ws.Range("I15").Value2 = 'Primera'
$ws.Range("J15").Value2 = 2960
$ws.Range("K15").Value2 = 2300
$ws.Range("L15").Value2 = 2500
$ws.Range("M15").Value2 = 2400
$ws.Range("N15").Value2 = '$/atado 1 a 1,5 kilos'
$ws.Range("O15").Value2 = 'Provincia del Elquí'
$ws.Range("P15").Value2 = 1600
$ws.Range("Q15").Value2 = 1.5
$ws.Range("R15").Value2 = 'Hortaliza'
